# Apply the "sliding window" update described by the diff:
#  - Drop the current first data row (row 2), shifting all existing rows up by one.
#  - Append 10 brand-new rows of data after the old last row, so that the sheet
#    ends up with data in rows 2-31 (30 data rows total).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to append, in order, as (x, y, z) triples.
$newRows = @(
    @(2.895918472953405, -3.924228730409128, 1.001562802687943),
    @(6.419859485349798, 3.571867339277037, 0.7840970496912497),
    @(-3.988556474879168, 3.858240864703046, 3.390691111053246),
    @(-2.792354390241082, 2.638935013093731, 1.279087337894736),
    @(-4.14442459281516, 9.648339557186993, -0.673063791604435),
    @(-1.790030563511108, 9.842405033572287, -8.527358793406183),
    @(6.688283160112809, -5.523453207983309, -0.2542450458244243),
    @(3.419822825325905, -5.20332591199646, -2.518984390917574),
    @(-0.0826715539042695, 0.4171818759705612, 2.113873891784312),
    @(-12.35876123801507, -14.20579128680002, 4.713338022646673),
    @(5.30543631401633, -16.49609409092704, 8.442643778335675)
)

# Delete the first data row (row 2), which shifts rows 3..21 up to become rows 2..20.
$ws.Rows.Item(2).Delete()

# After the delete, the old data occupies rows 2..20 (19 rows).
# Append the new rows starting at row 21 through row 31.
$startRow = 21
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
